$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.490.70"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").Value = "1.926.37"
$ws.Range("E3").Value = "  +1.72%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  +12.37%  "
$ws.Range("E6").Value = "  +5.00%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -0.63%  "
$ws.Range("E9").Value = "  +4.15%  "
$ws.Range("E10").Value = "  +4.56%  "
$ws.Range("E11").Value = "  +6.52%  "
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("D13").Value = "2.206.90"
$ws.Range("E13").Value = "  +1.77%  "
$ws.Range("E14").Value = "  +7.93%  "
$ws.Range("E15").Value = "  +4.56%  "
$ws.Range("D16").Value = "1.942.41"
$ws.Range("E16").Value = "  +2.84%  "
$ws.Range("E17").Value = "  +2.61%  "
$ws.Range("D18").Value = "35.493.50"
$ws.Range("E18").Value = "  +0.81%  "
$ws.Range("E19").Value = "  +5.05%  "
$ws.Range("E20").Value = "  +3.57%  "
$ws.Range("E21").Value = "  +1.75%  "
$ws.Range("E22").Value = "  +5.49%  "
$ws.Range("E23").Value = "  +8.24%  "
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("E25").Value = "  +1.74%  "
$ws.Range("E26").Value = "  +2.59%  "
$ws.Range("E27").Value = "  -0.82%  "
$ws.Range("E28").Value = "  +3.78%  "
$ws.Range("E29").Value = "  +6.99%  "
$ws.Range("E30").Value = "  +3.95%  "
$ws.Range("D31").Value = "4.125.93"
$ws.Range("E31").Value = "  +19.38%  "
$ws.Range("E32").Value = "  +6.66%  "
$ws.Range("E33").Value = "  +24.94%  "
$ws.Range("E34").Value = "  +14.34%  "
$ws.Range("E35").Value = "  +4.61%  "
$ws.Range("E36").Value = "  +3.65%  "
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("E38").Value = "  -1.90%  "
$ws.Range("E39").Value = "  +0.78%  "
$ws.Range("E40").Value = "  +10.25%  "
$ws.Range("E41").Value = "  +4.88%  "
$ws.Range("E42").Value = "  +9.92%  "
$ws.Range("E43").Value = "  +1.98%  "
$ws.Range("E44").Value = "  +2.12%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.352.22"
$ws.Range("E45").Value = "  +1.20%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E46").Value = "  +5.34%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("E47").Value = "  +5.45%  "
$ws.Range("B48").Value = "HuobiToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("E48").Value = "  +1.03%  "
$ws.Range("E49").Value = "  +0.59%  "
$ws.Range("E50").Value = "  -6.33%  "
$ws.Range("E51").Value = "  +7.28%  "

# Columns whose new value is a "clean" number-looking string (e.g. "41.08").
# Excel auto-converts such text to a real number on plain assignment, which
# would not match the source workbook (these are inlineStr text cells, not
# numeric cells). Force text entry via NumberFormat "@", then clear the
# temporary formatting so no extra style sticks to the cell (matches original,
# which has no explicit style on these cells).
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "255.10"
$ws.Range("D6").ClearFormats()
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.08"
$ws.Range("D8").ClearFormats()
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.359"
$ws.Range("D9").ClearFormats()
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.56"
$ws.Range("D10").ClearFormats()
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0753"
$ws.Range("D11").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "12.88"
$ws.Range("D14").ClearFormats()
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.722"
$ws.Range("D15").ClearFormats()
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.93"
$ws.Range("D17").ClearFormats()
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.65"
$ws.Range("D19").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "244.94"
$ws.Range("D21").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.06"
$ws.Range("D22").ClearFormats()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.11"
$ws.Range("D23").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.46"
$ws.Range("D25").ClearFormats()
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "168.32"
$ws.Range("D27").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.66"
$ws.Range("D28").ClearFormats()
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "18.88"
$ws.Range("D30").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.65"
$ws.Range("D33").ClearFormats()
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0585"
$ws.Range("D35").ClearFormats()
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.24"
$ws.Range("D36").ClearFormats()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.915"
$ws.Range("D38").ClearFormats()
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.03"
$ws.Range("D39").ClearFormats()
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.53"
$ws.Range("D40").ClearFormats()
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.66"
$ws.Range("D42").ClearFormats()
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0650"
$ws.Range("D44").ClearFormats()
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.47"
$ws.Range("D46").ClearFormats()
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.85"
$ws.Range("D47").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.43"
$ws.Range("D48").ClearFormats()
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.79"
$ws.Range("D49").ClearFormats()
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.38"
$ws.Range("D50").ClearFormats()
